# Update countries & provincias Spain
#
# The 'Pais' sheet lists countries with daily COVID-19 stats in
# columns B-H, ordered (roughly) by column B ('Casos totales')
# descending. This refresh updates several countries' figures;
# a few of them (Kenia, Tayikistan, Burundi) also change rank and
# need their row moved to keep the table in the right order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Sheet, $Row, $Values)
    $cols = @('A','B','C','D','E','F','G','H')
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $Sheet.Range("$($cols[$i])$Row").Value = $Values[$i]
    }
}

function Set-RowStats {
    # Writes Values into columns B..H of Row (leaves column A untouched)
    param($Sheet, $Row, $Values)
    $cols = @('B','C','D','E','F','G','H')
    for ($i = 0; $i -lt $Values.Count; $i++) {
        $Sheet.Range("$($cols[$i])$Row").Value = $Values[$i]
    }
}

# --- Countries whose figures changed but that keep their row ---
Set-RowStats $ws 4 @(1165953, 5179, 173910, 924448, 16475, 151, 67595)  # Estados Unidos
Set-RowStats $ws 15 @(57148, 434, 24416, 29126, 557, 40, 3606)  # Canada
Set-RowStats $ws 46 @(7764, 9, 3584, 3935, 62, 0, 245)  # Chequia
Set-RowStats $ws 58 @(4474, 179, 1936, 2075, 22, 4, 463)  # Argelia
Set-RowStats $ws 66 @(2626, 6, 1374, 1108, 37, 1, 144)  # Grecia
Set-RowStats $ws 72 @(2136, 18, 1319, 807, 8, 1, 10)  # Uzbekistan
Set-RowStats $ws 79 @(1618, 24, 308, 1237, 39, 1, 73)  # Bulgaria
Set-RowStats $ws 95 @(872, 8, 296, 561, 15, 0, 15)  # Republica de Chipre

# --- Rows whose ranking changed: rewrite rows 118-120, 146-164 and 198-199
#     in their final (post-refresh) order. ---
Set-RowValues $ws 118 @('Kenia', 465, 30, 167, 274, 2, 2, 24)
Set-RowValues $ws 119 @('Jamaica', 463, 31, 33, 422, 2, 0, 8)
Set-RowValues $ws 120 @('Jordania', 460, 0, 367, 84, 5, 0, 9)
Set-RowValues $ws 146 @('Tayikistan', 128, 52, 0, 126, 0, 0, 2)
Set-RowValues $ws 147 @('Guayana Francesa', 128, 0, 98, 29, 2, 0, 1)
Set-RowValues $ws 148 @('Togo', 124, 1, 67, 48, 0, 0, 9)
Set-RowValues $ws 149 @('Zambia', 124, 5, 78, 43, 1, 0, 3)
Set-RowValues $ws 150 @('Camboya', 122, 0, 120, 2, 1, 0, 0)
Set-RowValues $ws 151 @('Republica del Chad', 117, 0, 39, 68, 0, 0, 10)
Set-RowValues $ws 152 @('Trinidad yTobago', 116, 0, 88, 20, 0, 0, 8)
Set-RowValues $ws 153 @('Bermudas', 114, 0, 51, 56, 4, 0, 7)
Set-RowValues $ws 154 @('Suazilandia', 112, 4, 12, 99, 0, 0, 1)
Set-RowValues $ws 155 @('Aruba', 100, 0, 81, 17, 4, 0, 2)
Set-RowValues $ws 156 @('Monaco', 95, 0, 78, 13, 1, 0, 4)
Set-RowValues $ws 157 @('Benin', 90, 0, 42, 46, 0, 0, 2)
Set-RowValues $ws 158 @('Uganda', 88, 0, 52, 36, 0, 0, 0)
Set-RowValues $ws 159 @('Haiti', 85, 0, 10, 67, 0, 0, 8)
Set-RowValues $ws 160 @('Bahamas', 83, 0, 24, 48, 1, 0, 11)
Set-RowValues $ws 161 @('Guyana', 82, 0, 22, 51, 2, 0, 9)
Set-RowValues $ws 162 @('Liechtenstein', 82, 0, 55, 26, 0, 0, 1)
Set-RowValues $ws 163 @('Barbados', 81, 0, 44, 30, 4, 0, 7)
Set-RowValues $ws 164 @('Mozambique', 79, 0, 18, 61, 0, 0, 0)
Set-RowValues $ws 198 @('Burundi', 15, 0, 7, 7, 0, 0, 1)
Set-RowValues $ws 199 @('San Cristobal y Nieves', 15, 0, 8, 7, 0, 0, 0)
